$wb = $excel.ActiveWorkbook

$oldFile = "9a4d645e-ae4d-478e-b7df-4c0e2c368d03.md"
$newFile = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.md"
$newFile2 = "ffff7dd74055-cbb0-477a-ad16-7d3919fb31e3.md"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bce338d628024515be7284f0f35e0c26743c5b8e/e2e/"

$oldXlfZh = "9a4d645e-ae4d-478e-b7df-4c0e2c368d03.953250b68628d052239c54f79af3d2fee25e43e3.zh-cn.xlf"
$newXlfZh = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.1f5c1ef0c073e683166b6b5c438544bd79b4898e.zh-cn.xlf"
$oldXlfDe = "9a4d645e-ae4d-478e-b7df-4c0e2c368d03.953250b68628d052239c54f79af3d2fee25e43e3.de-de.xlf"
$newXlfDe = "4066dd2f-145f-4bc4-9a29-b8822e9b16ac.1f5c1ef0c073e683166b6b5c438544bd79b4898e.de-de.xlf"

# =========================================================================
# Sheet "Overview"
# =========================================================================
$ws = $wb.Worksheets.Item("Overview")

# -- update existing row 2 (renamed source file + new handoff datetime) --
$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = "e2e\" + $newFile
$ws.Range("G2").Value = "2017-02-09 16:04:58"

$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), ($baseUrl + $newFile), "", "", ("e2e\" + $newFile)) | Out-Null

# -- append new row 3 for the newly handed-off file --
$ws.Range("A3").Value = $newFile2
$ws.Range("B3").Value = "e2e\" + $newFile2
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2017-02-09 16:04:58"

$ws.Hyperlinks.Add($ws.Range("B3"), ($baseUrl + $newFile2), "", "", ("e2e\" + $newFile2)) | Out-Null

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

# -- update existing row 2 --
$ws.Range("A2").Value = $newFile
$ws.Range("G2").Value = $newXlfZh
$ws.Range("H2").Value = "2017-02-09 16:04:41"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), ($baseUrl + $newFile), "", "", $newFile) | Out-Null

# -- append new row 3 --
$ws.Range("A3").Value = $newFile2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = $newXlfZh
$ws.Range("H3").Value = "2017-02-09 16:04:41"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "0001-01-01 00:00:00"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "True"
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = "False"
$ws.Range("R3").Value = ""

$ws.Hyperlinks.Add($ws.Range("A3"), ($baseUrl + $newFile2), "", "", $newFile2) | Out-Null

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:R3"))

# =========================================================================
# Sheet "de-de"
# =========================================================================
$ws = $wb.Worksheets.Item("de-de")

# -- update existing row 2 --
$ws.Range("A2").Value = $newFile
$ws.Range("G2").Value = $newXlfDe
$ws.Range("H2").Value = "2017-02-09 16:04:58"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), ($baseUrl + $newFile), "", "", $newFile) | Out-Null

# -- append new row 3 --
$ws.Range("A3").Value = $newFile2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = $newXlfDe
$ws.Range("H3").Value = "2017-02-09 16:04:58"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "0001-01-01 00:00:00"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "True"
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = "False"
$ws.Range("R3").Value = ""

$ws.Hyperlinks.Add($ws.Range("A3"), ($baseUrl + $newFile2), "", "", $newFile2) | Out-Null

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:R3"))

Write-Host "Report regenerated for handoff."
